$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Seed the new shared strings in the same order the original authoring
#     produced them (cosmetic, but keeps sharedStrings.xml close to target) ---
$ws.Range("E10").Value = "Tegshig"
$ws.Range("E12").Value = "Tim"
$ws.Range("D12").Value = "Start design of Node-Red"
$ws.Range("E11").Value = "Luca"
$ws.Range("D28").Value = "Finalize the website"
$ws.Range("A10").Value = "Discuss the project"
$ws.Range("D44").Value = "Optimizing website"
$ws.Range("B68").Value = "Upload Project"
$ws.Range("B5").Value = "K"

# --- Week 1 (rows 9-12) ---
$ws.Range("B10").Value = "grupo"
$ws.Range("C10").Value = 1
$ws.Range("F10").Value = 2.25
$ws.Range("F11").Value = 2.25
$ws.Range("F12").Value = 2.25

# --- Week 2 (rows 25-28) ---
$ws.Range("C26").Value = 2
$ws.Range("E26").Value = "Tegshig"
$ws.Range("F26").Value = 2.25

$ws.Range("A27").Value = "Discuss the project"
$ws.Range("B27").Value = "grupo"
$ws.Range("C27").Value = 1
$ws.Range("E27").Value = "Luca"
$ws.Range("F27").Value = 2.25

$ws.Range("E28").Value = "Tim"
$ws.Range("F28").Value = 2.25

# --- Week 3 (rows 41-44) ---
$ws.Range("A42").Value = "Discuss the project"
$ws.Range("B42").Value = "grupo"
$ws.Range("C42").Value = 1
$ws.Range("D42").Value = "Trouble shooting"
$ws.Range("E42").Value = "Tegshig"
$ws.Range("F42").Value = 2.25

$ws.Range("D43").Value = "Finalising project"
$ws.Range("E43").Value = "Luca"
$ws.Range("F43").Value = 2.25

$ws.Range("E44").Value = "Tim"
$ws.Range("F44").Value = 2.25

# --- Week 4 (rows 57-60) ---
$ws.Range("B58").Value = "grupo"
$ws.Range("C58").Value = 1

$ws.Range("B59").Value = "grupo"
$ws.Range("C59").Value = 1

$ws.Range("A60").Value = "Discuss the project"
$ws.Range("B60").Value = "grupo"
$ws.Range("C60").Value = 1

# --- Formatting ---
$ws.Range("C8").WrapText = $true
$ws.Range("F8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 35.25

$ws.Columns.Item(3).ColumnWidth = 24.307291666666668
$ws.Columns.Item(6).ColumnWidth = 25.166666666666668

# --- Selection ---
$ws.Range("B5").Select()
